$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# A95: date-time value, reuse the existing date style (same as column A in prior rows)
$ws.Cells.Item($row, 1).Value = 45446.2916666667
$ws.Cells.Item(94, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# B95-F95: plain numeric values
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 1

# G95: text "1" (looks numeric, force text storage without leaving a custom number format behind)
$ws.Cells.Item($row, 7).Value = "'1"
$ws.Cells.Item($row, 7).Style = "Normal"

# H95: ticker text
$ws.Cells.Item($row, 8).Value = "YKY.MI"
